# Edit script: add new survey wave column (11.-17. 10. 2021) to both sheets
# and bump the "aktualizace" (updated) footer date from 6. 10. 2021 to 20. 10. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "data" -- new column AJ (36th column), rows 1-301 + footer row 302
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Header cell AJ1: same look as the existing header row (bold, centered,
# top-aligned, thin border) -- matches AI1's style.
$hdr1 = $ws1.Range("AJ1")
$hdr1.Font.Bold = $true
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4160
$hdr1.Borders.LineStyle = 1
$hdr1.Value2 = "11.–17. 10. 2021"

# Data values AJ2:AJ301, in row order.
$sheet1Values = @(0.8100000000000001,0.38,0.3,0.24,0.05,0.41,0.11,0.13,0.09,0.05,0.64,0.28,0.13,0.12,0.03,0.78,0.85,0.78,0.82,0.78,0.84,0.8100000000000001,0.87,0.79,0.82,0.79,0.8100000000000001,0.8,0.84,0.8,0.87,0.86,0.9399999999999999,0.72,0.5,0.36,0.33,0.4,0.39,0.35,0.47,0.37,0.36,0.38,0.38,0.27,0.35,0.6,0.68,0.4,0.42,0.61,0.42,0.26,0.27,0.35,0.24,0.29,0.37,0.35,0.28,0.29,0.25,0.35,0.3,0.28,0.31,0.32,0.24,0.23,0.22,0.54,0.31,0.21,0.22,0.26,0.22,0.25,0.36,0.17,0.23,0.23,0.24,0.24,0.2,0.26,0.3,0.25,0.21,0.28,0.26,0.07000000000000001,0.05,0.03,0.05,0.04,0.05,0.08,0.05,0.04,0.06,0.04,0.03,0.05,0.04,0.13,0.05,0.12,0.06,0.04,0.55,0.38,0.34,0.51,0.38,0.35,0.49,0.41,0.38,0.49,0.32,0.36,0.4,0.5,0.52,0.46,0.58,0.54,0.32,0.19,0.12,0.05,0.16,0.09,0.08,0.18,0.08,0.1,0.12,0.09,0.09,0.11,0.14,0.14,0.14,0.11,0.23,0.03,0.15,0.13,0.11,0.14,0.11,0.14,0.18,0.13,0.11,0.14,0.11,0.13,0.13,0.15,0.1,0.15,0.04,0.15,0.1,0.12,0.1,0.06,0.11,0.08,0.08,0.12,0.08,0.08,0.07000000000000001,0.1,0.1,0.09,0.07000000000000001,0.06,0.12,0.1,0.015,0.07000000000000001,0.08,0.04,0.03,0.06,0.04,0.05,0.08,0.03,0.04,0.05,0.04,0.04,0.05,0.07000000000000001,0.04,0.06,0.01,0.08,0.06,0.6899999999999999,0.61,0.62,0.66,0.61,0.65,0.67,0.68,0.61,0.62,0.65,0.6899999999999999,0.5600000000000001,0.61,0.58,0.64,0.63,0.63,0.5,0.37,0.28,0.22,0.33,0.24,0.3,0.37,0.28,0.26,0.27,0.29,0.28,0.28,0.25,0.31,0.32,0.26,0.24,0.28,0.13,0.14,0.12,0.12,0.12,0.15,0.21,0.11,0.12,0.15,0.11,0.11,0.11,0.17,0.21,0.17,0.09,0.15,0.09,0.19,0.14,0.07000000000000001,0.15,0.11,0.11,0.19,0.1,0.11,0.13,0.12,0.12,0.1,0.15,0.17,0.15,0.13,0.06,0.19,0.06,0.04,0.015,0.04,0.015,0.05,0.05,0.03,0.03,0.03,0.03,0.04,0.015,0.03,0.015,0.04,0.04,0.05,0.06)
for ($i = 0; $i -lt $sheet1Values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 36).Value2 = $sheet1Values[$i]
}

# Footer label (row 302, column A): bump the update date.
$footer1 = $ws1.Range("A302")
$footer1.Value2 = $footer1.Value2.Replace("6. 10. 2021", "20. 10. 2021")

# ---------------------------------------------------------------------------
# Sheet 2: "pocetR" -- new column AI (35th column), rows 1-21 + footer row 22
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Header cell AI1: same look as the existing header row.
$hdr2 = $ws2.Range("AI1")
$hdr2.Font.Bold = $true
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4160
$hdr2.Borders.LineStyle = 1
$hdr2.Value2 = "11.–17. 10. 2021"

# Data values AI2:AI21, in row order.
$sheet2Values = @(1836,454,670,712,522,819,495,296,316,1224,895,941,960,418,216,242,809,70,75,87)
for ($i = 0; $i -lt $sheet2Values.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 35).Value2 = $sheet2Values[$i]
}

# Trailing blank placeholder cell AI22, matching the other blank cells
# already present across row 22 (B22:AH22).
$ws2.Range("AI22").Style = "Normal"

# Footer label (row 22, column A): bump the update date.
$footer2 = $ws2.Range("A22")
$footer2.Value2 = $footer2.Value2.Replace("6. 10. 2021", "20. 10. 2021")
